$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gnai2"
$ws.Cells.Item(2,3).Value = "Oprm1"
$ws.Cells.Item(2,4).Value = "M1"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 195.0792385
$ws.Cells.Item(2,8).Value = 390.158477
$ws.Cells.Item(2,9).Value = 0.2640605522989327
$ws.Cells.Item(2,10).Value = 0.1982306263353075
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.05215766666666666
$ws.Cells.Item(2,14).Value = 0.156473
$ws.Cells.Item(2,15).Value = 0.1010355835763341
$ws.Cells.Item(2,16).Value = 0.1010355835763341
$ws.Cells.Item(2,17).Value = 10.17487789527017
$ws.Cells.Item(2,18).Value = 61.049267371621
$ws.Cells.Item(2,19).Value = 0.02667951200101175
$ws.Cells.Item(2,20).Value = 0.02002834701449002

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gnai2"
$ws.Cells.Item(3,3).Value = "Oprm1"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 195.0792385
$ws.Cells.Item(3,8).Value = 390.158477
$ws.Cells.Item(3,9).Value = 0.2640605522989327
$ws.Cells.Item(3,10).Value = 0.1982306263353075
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.464073
$ws.Cells.Item(3,14).Value = 1.392219
$ws.Cells.Item(3,15).Value = 0.8989644164236659
$ws.Cells.Item(3,16).Value = 0.8989644164236659
$ws.Cells.Item(3,17).Value = 90.5310074484105
$ws.Cells.Item(3,18).Value = 543.1860446904631
$ws.Cells.Item(3,19).Value = 0.2373810402979209
$ws.Cells.Item(3,20).Value = 0.1782022793208175

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gnai2"
$ws.Cells.Item(4,3).Value = "Oprm1"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 62.40792233333334
$ws.Cells.Item(4,8).Value = 187.223767
$ws.Cells.Item(4,9).Value = 0.08447577797556809
$ws.Cells.Item(4,10).Value = 0.09512412720758515
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.05215766666666666
$ws.Cells.Item(4,14).Value = 0.156473
$ws.Cells.Item(4,15).Value = 0.1010355835763341
$ws.Cells.Item(4,16).Value = 0.1010355835763341
$ws.Cells.Item(4,17).Value = 3.255051610421222
$ws.Cells.Item(4,18).Value = 29.295464493791
$ws.Cells.Item(4,19).Value = 0.00853505952582635
$ws.Cells.Item(4,20).Value = 0.009610921704607804

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gnai2"
$ws.Cells.Item(5,3).Value = "Oprm1"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 62.40792233333334
$ws.Cells.Item(5,8).Value = 187.223767
$ws.Cells.Item(5,9).Value = 0.08447577797556809
$ws.Cells.Item(5,10).Value = 0.09512412720758515
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.464073
$ws.Cells.Item(5,14).Value = 1.392219
$ws.Cells.Item(5,15).Value = 0.8989644164236659
$ws.Cells.Item(5,16).Value = 0.8989644164236659
$ws.Cells.Item(5,17).Value = 28.961831740997
$ws.Cells.Item(5,18).Value = 260.656485668973
$ws.Cells.Item(5,19).Value = 0.07594071844974173
$ws.Cells.Item(5,20).Value = 0.08551320550297734

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Gnai2"
$ws.Cells.Item(6,3).Value = "Oprm1"
$ws.Cells.Item(6,4).Value = "M1"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 142.8621113333333
$ws.Cells.Item(6,8).Value = 428.586334
$ws.Cells.Item(6,9).Value = 0.1933791023142199
$ws.Cells.Item(6,10).Value = 0.2177549443006804
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.05215766666666666
$ws.Cells.Item(6,14).Value = 0.156473
$ws.Cells.Item(6,15).Value = 0.1010355835763341
$ws.Cells.Item(6,16).Value = 0.1010355835763341
$ws.Cells.Item(6,17).Value = 7.451354382220222
$ws.Cells.Item(6,18).Value = 67.06218943998199
$ws.Cells.Item(6,19).Value = 0.01953817045378482
$ws.Cells.Item(6,20).Value = 0.02200099787405137

# Row 7
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Gnai2"
$ws.Cells.Item(7,3).Value = "Oprm1"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 142.8621113333333
$ws.Cells.Item(7,8).Value = 428.586334
$ws.Cells.Item(7,9).Value = 0.1933791023142199
$ws.Cells.Item(7,10).Value = 0.2177549443006804
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.464073
$ws.Cells.Item(7,14).Value = 1.392219
$ws.Cells.Item(7,15).Value = 0.8989644164236659
$ws.Cells.Item(7,16).Value = 0.8989644164236659
$ws.Cells.Item(7,17).Value = 66.298448592794
$ws.Cells.Item(7,18).Value = 596.686037335146
$ws.Cells.Item(7,19).Value = 0.1738409318604351
$ws.Cells.Item(7,20).Value = 0.195753946426629

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Gnai2"
$ws.Cells.Item(8,3).Value = "Oprm1"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 141.6168416666667
$ws.Cells.Item(8,8).Value = 424.850525
$ws.Cells.Item(8,9).Value = 0.1916934970264942
$ws.Cells.Item(8,10).Value = 0.2158568649262854
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.05215766666666666
$ws.Cells.Item(8,14).Value = 0.156473
$ws.Cells.Item(8,15).Value = 0.1010355835763341
$ws.Cells.Item(8,16).Value = 0.1010355835763341
$ws.Cells.Item(8,17).Value = 7.38640402203611
$ws.Cells.Item(8,18).Value = 66.477636198325
$ws.Cells.Item(8,19).Value = 0.0193678643398601
$ws.Cells.Item(8,20).Value = 0.02180922431678516

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Gnai2"
$ws.Cells.Item(9,3).Value = "Oprm1"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 141.6168416666667
$ws.Cells.Item(9,8).Value = 424.850525
$ws.Cells.Item(9,9).Value = 0.1916934970264942
$ws.Cells.Item(9,10).Value = 0.2158568649262854
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.464073
$ws.Cells.Item(9,14).Value = 1.392219
$ws.Cells.Item(9,15).Value = 0.8989644164236659
$ws.Cells.Item(9,16).Value = 0.8989644164236659
$ws.Cells.Item(9,17).Value = 65.720552562775
$ws.Cells.Item(9,18).Value = 591.4849730649751
$ws.Cells.Item(9,19).Value = 0.1723256326866341
$ws.Cells.Item(9,20).Value = 0.1940476406095002

# Row 10
$ws.Cells.Item(10,1).Value = "Neutro"
$ws.Cells.Item(10,2).Value = "Gnai2"
$ws.Cells.Item(10,3).Value = "Oprm1"
$ws.Cells.Item(10,4).Value = "M1"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 143.783834
$ws.Cells.Item(10,8).Value = 431.351502
$ws.Cells.Item(10,9).Value = 0.1946267522348261
$ws.Cells.Item(10,10).Value = 0.2191598631141254
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.05215766666666666
$ws.Cells.Item(10,14).Value = 0.156473
$ws.Cells.Item(10,15).Value = 0.1010355835763341
$ws.Cells.Item(10,16).Value = 0.1010355835763341
$ws.Cells.Item(10,17).Value = 7.499429285827333
$ws.Cells.Item(10,18).Value = 67.494863572446
$ws.Cells.Item(10,19).Value = 0.01966422749161224
$ws.Cells.Item(10,20).Value = 0.02214294466624515

# Row 11
$ws.Cells.Item(11,1).Value = "Neutro"
$ws.Cells.Item(11,2).Value = "Gnai2"
$ws.Cells.Item(11,3).Value = "Oprm1"
$ws.Cells.Item(11,4).Value = "M2"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 143.783834
$ws.Cells.Item(11,8).Value = 431.351502
$ws.Cells.Item(11,9).Value = 0.1946267522348261
$ws.Cells.Item(11,10).Value = 0.2191598631141254
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.464073
$ws.Cells.Item(11,14).Value = 1.392219
$ws.Cells.Item(11,15).Value = 0.8989644164236659
$ws.Cells.Item(11,16).Value = 0.8989644164236659
$ws.Cells.Item(11,17).Value = 66.726195195882
$ws.Cells.Item(11,18).Value = 600.535756762938
$ws.Cells.Item(11,19).Value = 0.1749625247432138
$ws.Cells.Item(11,20).Value = 0.1970169184478802

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Gnai2"
$ws.Cells.Item(12,3).Value = "Oprm1"
$ws.Cells.Item(12,4).Value = "M1"
$ws.Cells.Item(12,5).Value = 2
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 53.01711450000001
$ws.Cells.Item(12,8).Value = 106.034229
$ws.Cells.Item(12,9).Value = 0.07176431814995911
$ws.Cells.Item(12,10).Value = 0.05387357411601602
$ws.Cells.Item(12,11).Value = 1
$ws.Cells.Item(12,12).Value = 0.3333333333333333
$ws.Cells.Item(12,13).Value = 0.05215766666666666
$ws.Cells.Item(12,14).Value = 0.156473
$ws.Cells.Item(12,15).Value = 0.1010355835763341
$ws.Cells.Item(12,16).Value = 0.1010355835763341
$ws.Cells.Item(12,17).Value = 2.7652489857195
$ws.Cells.Item(12,18).Value = 16.591493914317
$ws.Cells.Item(12,19).Value = 0.007250749764238822
$ws.Cells.Item(12,20).Value = 0.005443148000154565

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Gnai2"
$ws.Cells.Item(13,3).Value = "Oprm1"
$ws.Cells.Item(13,4).Value = "M2"
$ws.Cells.Item(13,5).Value = 2
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 53.01711450000001
$ws.Cells.Item(13,8).Value = 106.034229
$ws.Cells.Item(13,9).Value = 0.07176431814995911
$ws.Cells.Item(13,10).Value = 0.05387357411601602
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.464073
$ws.Cells.Item(13,14).Value = 1.392219
$ws.Cells.Item(13,15).Value = 0.8989644164236659
$ws.Cells.Item(13,16).Value = 0.8989644164236659
$ws.Cells.Item(13,17).Value = 24.6038113773585
$ws.Cells.Item(13,18).Value = 147.622868264151
$ws.Cells.Item(13,19).Value = 0.06451356838572028
$ws.Cells.Item(13,20).Value = 0.04843042611586145
